$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new question about Jupyter QC notebooks
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "How comfortable are you using Jupyter QC notebooks to review the quality of MS/MS identifications"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1

# Row 8: new question about 'git pull'
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "How comfortable are you about using the 'git pull' command to update the IBIP21 folder"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 2

# Row 9: reuses existing question text (shared string index 3)
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "If a colleague is interested in a mutation in a given gene, how confident do you feel about changing the protein sequence to find the mutated peptide"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1

# Row 10: new question about Ensembl
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "If a student asks you how to find the consequence of a variant on a protein, how confident do you feel about explaining them how to find it in Ensembl?"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 3

$ws.Range("E11").Select()
